$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.816.03"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "'2.247.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'112.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "
$ws.Range("D6").Value = "'295.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.53%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").Value = "'44.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.98%  "
$ws.Range("D11").Value = "'0.0925"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "'54.37"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "'9.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.30%  "
$ws.Range("E14").Value = "  +20.75%  "
$ws.Range("E15").Value = "  -1.08%  "
$ws.Range("D16").Value = "'15.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.17%  "
$ws.Range("D17").Value = "'2.585.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "'2.249.33"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.51%  "
$ws.Range("D19").Value = "'42.748.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.0000106"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.00%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'7.22"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.38%  "
$ws.Range("D22").Value = "'74.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.20%  "
$ws.Range("D23").Value = "'3.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.22%  "
$ws.Range("D24").Value = "'2.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.99%  "
$ws.Range("D25").Value = "'250.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.00%  "
$ws.Range("E26").Value = "  -3.66%  "
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "'11.53"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.95%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'37.77"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.31%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'22.21"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.23%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'174.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D32").Value = "'2.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.89%  "
$ws.Range("D33").Value = "'3.16"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.50%  "
$ws.Range("D34").Value = "'0.0889"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("D35").Value = "'5.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("E36").Value = "  +9.10%  "
$ws.Range("D37").Value = "'4.27"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.30%  "
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").Value = "'0.0378"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.14%  "
$ws.Range("D40").Value = "'0.104"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.65%  "
$ws.Range("D41").Value = "'2.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.08%  "
$ws.Range("D42").Value = "'72.05"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("D43").Value = "'0.231"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "'12.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.32%  "
$ws.Range("D46").Value = "'1.32"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.03%  "
$ws.Range("D47").Value = "'5.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.11%  "
$ws.Range("E48").Value = "  +2.15%  "
$ws.Range("D49").Value = "'105.26"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.46%  "
$ws.Range("D50").Value = "'8.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.81%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.34%  "
